# Auto-generated: rebuild sheet1 data grid per target diff (adds parent.key / parent.typeId columns
# and reorders/extends category rows to include parent hierarchy + new categories).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 21,8
$data[0,0] = 'data-object'
$data[0,1] = 'key'
$data[0,2] = 'description.en-US'
$data[0,3] = 'externalId'
$data[0,4] = 'name.en-US'
$data[0,5] = 'slug.en-US'
$data[0,6] = 'parent.key'
$data[0,7] = 'parent.typeId'
$data[1,0] = 'category'
$data[1,1] = 'abilityLevelKey'
$data[1,2] = 'abilityLevelDescription'
$data[1,3] = 'abilityLevelId'
$data[1,4] = 'AbilityLevel'
$data[1,5] = 'abilityLevelSlug'
$data[1,6] = ''
$data[1,7] = ''
$data[2,0] = 'category'
$data[2,1] = 'brandKey'
$data[2,2] = 'brandDescription'
$data[2,3] = 'brandId'
$data[2,4] = 'Brand'
$data[2,5] = 'brandSlug'
$data[2,6] = ''
$data[2,7] = ''
$data[3,0] = 'category'
$data[3,1] = 'ageGroupKey'
$data[3,2] = 'ageGroupDescription'
$data[3,3] = 'ageGroupId'
$data[3,4] = 'AgeGroup'
$data[3,5] = 'ageGroupSlug'
$data[3,6] = ''
$data[3,7] = ''
$data[4,0] = 'category'
$data[4,1] = 'terrainKey'
$data[4,2] = 'terrainDescription'
$data[4,3] = 'terrainId'
$data[4,4] = 'Terrain'
$data[4,5] = 'terrainSlug'
$data[4,6] = ''
$data[4,7] = ''
$data[5,0] = 'category'
$data[5,1] = 'Advanced-ExpertKey'
$data[5,2] = 'Advanced-ExpertDescription'
$data[5,3] = 'Advanced-ExpertId'
$data[5,4] = 'Advanced-Expert'
$data[5,5] = 'Advanced-ExpertSlug'
$data[5,6] = 'abilityLevelKey'
$data[5,7] = 'category'
$data[6,0] = 'category'
$data[6,1] = 'LibTechKey'
$data[6,2] = 'LibTechDescription'
$data[6,3] = 'LibTechId'
$data[6,4] = 'LibTech'
$data[6,5] = 'LibTechSlug'
$data[6,6] = 'brandKey'
$data[6,7] = 'category'
$data[7,0] = 'category'
$data[7,1] = 'adult-maleKey'
$data[7,2] = 'adult-maleDescription'
$data[7,3] = 'adult-maleId'
$data[7,4] = 'Adult-male'
$data[7,5] = 'adult-maleSlug'
$data[7,6] = 'ageGroupKey'
$data[7,7] = 'category'
$data[8,0] = 'category'
$data[8,1] = 'FreestyleKey'
$data[8,2] = 'FreestyleDescription'
$data[8,3] = 'FreestyleId'
$data[8,4] = 'Freestyle'
$data[8,5] = 'FreestyleSlug'
$data[8,6] = 'terrainKey'
$data[8,7] = 'category'
$data[9,0] = 'category'
$data[9,1] = 'All-MountainKey'
$data[9,2] = 'All-MountainDescription'
$data[9,3] = 'All-MountainId'
$data[9,4] = 'All-Mountain'
$data[9,5] = 'All-MountainSlug'
$data[9,6] = 'terrainKey'
$data[9,7] = 'category'
$data[10,0] = 'category'
$data[10,1] = 'FreerideKey'
$data[10,2] = 'FreerideDescription'
$data[10,3] = 'FreerideId'
$data[10,4] = 'Freeride'
$data[10,5] = 'FreerideSlug'
$data[10,6] = 'terrainKey'
$data[10,7] = 'category'
$data[11,0] = 'category'
$data[11,1] = 'Intermediate-AdvancedKey'
$data[11,2] = 'Intermediate-AdvancedDescription'
$data[11,3] = 'Intermediate-AdvancedId'
$data[11,4] = 'Intermediate-Advanced'
$data[11,5] = 'Intermediate-AdvancedSlug'
$data[11,6] = 'abilityLevelKey'
$data[11,7] = 'category'
$data[12,0] = 'category'
$data[12,1] = 'RideKey'
$data[12,2] = 'RideDescription'
$data[12,3] = 'RideId'
$data[12,4] = 'Ride'
$data[12,5] = 'RideSlug'
$data[12,6] = 'brandKey'
$data[12,7] = 'category'
$data[13,0] = 'category'
$data[13,1] = 'adult-anyKey'
$data[13,2] = 'adult-anyDescription'
$data[13,3] = 'adult-anyId'
$data[13,4] = 'Adult-any'
$data[13,5] = 'adult-anySlug'
$data[13,6] = 'ageGroupKey'
$data[13,7] = 'category'
$data[14,0] = 'category'
$data[14,1] = 'CAPiTAKey'
$data[14,2] = 'CAPiTADescription'
$data[14,3] = 'CAPiTAId'
$data[14,4] = 'CAPiTA'
$data[14,5] = 'CAPiTASlug'
$data[14,6] = 'brandKey'
$data[14,7] = 'category'
$data[15,0] = 'category'
$data[15,1] = 'SeasonKey'
$data[15,2] = 'SeasonDescription'
$data[15,3] = 'SeasonId'
$data[15,4] = 'Season'
$data[15,5] = 'SeasonSlug'
$data[15,6] = 'brandKey'
$data[15,7] = 'category'
$data[16,0] = 'category'
$data[16,1] = 'PowderKey'
$data[16,2] = 'PowderDescription'
$data[16,3] = 'PowderId'
$data[16,4] = 'Powder'
$data[16,5] = 'PowderSlug'
$data[16,6] = 'terrainKey'
$data[16,7] = 'category'
$data[17,0] = 'category'
$data[17,1] = 'Beginner-IntermediateKey'
$data[17,2] = 'Beginner-IntermediateDescription'
$data[17,3] = 'Beginner-IntermediateId'
$data[17,4] = 'Beginner-Intermediate'
$data[17,5] = 'Beginner-IntermediateSlug'
$data[17,6] = 'abilityLevelKey'
$data[17,7] = 'category'
$data[18,0] = 'category'
$data[18,1] = 'RossignolKey'
$data[18,2] = 'RossignolDescription'
$data[18,3] = 'RossignolId'
$data[18,4] = 'Rossignol'
$data[18,5] = 'RossignolSlug'
$data[18,6] = 'brandKey'
$data[18,7] = 'category'
$data[19,0] = 'category'
$data[19,1] = 'adult-femaleKey'
$data[19,2] = 'adult-femaleDescription'
$data[19,3] = 'adult-femaleId'
$data[19,4] = 'Adult-female'
$data[19,5] = 'adult-femaleSlug'
$data[19,6] = 'ageGroupKey'
$data[19,7] = 'category'
$data[20,0] = 'category'
$data[20,1] = 'K2Key'
$data[20,2] = 'K2Description'
$data[20,3] = 'K2Id'
$data[20,4] = 'K2'
$data[20,5] = 'K2Slug'
$data[20,6] = 'brandKey'
$data[20,7] = 'category'

$ws.Range("A1:H21").Value = $data

